$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-03 06:42:11"
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-03 06:42:05"
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-03 06:42:11"
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
